# cryptos list refresh (GitHub Actions bot)
# Updates Price (D) / Volume(1h) (E) for each coin row, and the two
# rows (49/50) that swapped rank order (RenderToken <-> RocketPoolETH).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.238.54'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '1.858.15'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''0.7027'
$ws.Range("E5").Value = '  +2.33%  '
$ws.Range("D6").Value = '''238.07'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''0.08016'
$ws.Range("E8").Value = '  +4.17%  '
$ws.Range("D9").Value = '''0.3022'
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("D10").Value = '''23.55'
$ws.Range("E10").Value = '  +1.88%  '
$ws.Range("D11").Value = '''0.08185'
$ws.Range("E11").Value = '  +0.43%  '
$ws.Range("D12").Value = '1.904.24'
$ws.Range("E12").Value = '  +2.88%  '
$ws.Range("D13").Value = '''5.199'
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").Value = '''0.7077'
$ws.Range("E14").Value = '  -1.88%  '
$ws.Range("D15").Value = '''89.60'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("D16").Value = '29.440.43'
$ws.Range("E16").Value = '  +0.99%  '
$ws.Range("D17").Value = '''5.833'
$ws.Range("E17").Value = '  +1.92%  '
$ws.Range("D18").Value = '''0.000007900'
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("D19").Value = '''13.28'
$ws.Range("E19").Value = '  +0.94%  '
$ws.Range("D20").Value = '''238.12'
$ws.Range("E20").Value = '  +1.85%  '
$ws.Range("D21").Value = '2.180.47'
$ws.Range("E21").Value = '  +4.06%  '
$ws.Range("D22").Value = '''0.9994'
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = '''1.001'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '''7.477'
$ws.Range("E24").Value = '  -0.33%  '
$ws.Range("D25").Value = '''162.80'
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("D26").Value = '''8.891'
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("D27").Value = '''0.1436'
$ws.Range("E27").Value = '  +1.22%  '
$ws.Range("D28").Value = '''18.11'
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").Value = '''1.921'
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("D30").Value = '''1.415'
$ws.Range("E30").Value = '  +1.08%  '
$ws.Range("D31").Value = '''1.475'
$ws.Range("E31").Value = '  -0.61%  '
$ws.Range("D32").Value = '''4.376'
$ws.Range("E32").Value = '  -2.87%  '
$ws.Range("D33").Value = '''4.026'
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").Value = '''0.05200'
$ws.Range("E34").Value = '  +0.37%  '
$ws.Range("D35").Value = '''1.163'
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("D36").Value = '''0.7183'
$ws.Range("E36").Value = '  +2.49%  '
$ws.Range("D37").Value = '''1.003'
$ws.Range("E37").Value = '  -1.70%  '
$ws.Range("D38").Value = '''2.699'
$ws.Range("E38").Value = '  +1.82%  '
$ws.Range("D39").Value = '''0.01852'
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("D40").Value = '''2.726'
$ws.Range("E40").Value = '  +1.72%  '
$ws.Range("D41").Value = '''0.9450'
$ws.Range("E41").Value = '  +3.77%  '
$ws.Range("D42").Value = '1.152.65'
$ws.Range("E42").Value = '  +5.38%  '
$ws.Range("E43").Value = '  +0.64%  '
$ws.Range("D44").Value = '''0.4268'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").Value = '''70.88'
$ws.Range("E45").Value = '  +0.89%  '
$ws.Range("D46").Value = '''0.9999'
$ws.Range("D47").Value = '''103.01'
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").Value = '''0.5304'
$ws.Range("E48").Value = '  -3.93%  '
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.042.17'
$ws.Range("E49").Value = '  +2.22%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '''1.761'
$ws.Range("E50").Value = '  +0.61%  '
$ws.Range("D51").Value = '''9.179'
$ws.Range("E51").Value = '  +0.64%  '
